$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-10-28 06:37:23"

for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
